$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force Text format first so the literal string (with its trailing zeros /
# separators) is preserved, exactly like the original sheet held plain text.
$numericRiskCells = @('D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D13', 'D14', 'D15', 'D16', 'D18', 'D19', 'D21', 'D22', 'D24', 'D25', 'D26', 'D27', 'D28', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $numericRiskCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '29.462.05'
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').Value = '1.872.06'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('D4').Value = '0.9995'
$ws.Range('D5').Value = '0.7083'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('D6').Value = '243.91'
$ws.Range('E6').Value = '  +0.57%  '
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').Value = '0.3165'
$ws.Range('E8').Value = '  +0.73%  '
$ws.Range('D9').Value = '0.07895'
$ws.Range('E9').Value = '  -1.65%  '
$ws.Range('D10').Value = '24.66'
$ws.Range('E10').Value = '  -1.87%  '
$ws.Range('D11').Value = '0.08002'
$ws.Range('E11').Value = '  -3.91%  '
$ws.Range('D12').Value = '1.882.79'
$ws.Range('E12').Value = '  -0.40%  '
$ws.Range('D13').Value = '5.229'
$ws.Range('E13').Value = '  -0.77%  '
$ws.Range('D14').Value = '94.23'
$ws.Range('D15').Value = '0.7061'
$ws.Range('E15').Value = '  -1.71%  '
$ws.Range('D16').Value = '6.523'
$ws.Range('E16').Value = '  +2.65%  '
$ws.Range('D17').Value = '29.470.74'
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('D18').Value = '0.000008368'
$ws.Range('E18').Value = '  -4.04%  '
$ws.Range('D19').Value = '257.68'
$ws.Range('E19').Value = '  +6.00%  '
$ws.Range('D20').Value = '2.127.81'
$ws.Range('E20').Value = '  -1.06%  '
$ws.Range('D21').Value = '13.24'
$ws.Range('E21').Value = '  -0.58%  '
$ws.Range('D22').Value = '0.9997'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('E23').Value = '  -2.89%  '
$ws.Range('D24').Value = '0.9993'
$ws.Range('E24').Value = '  -0.22%  '
$ws.Range('D25').Value = '0.1565'
$ws.Range('E25').Value = '  -0.60%  '
$ws.Range('D26').Value = '9.091'
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').Value = '160.93'
$ws.Range('E27').Value = '  -1.70%  '
$ws.Range('D28').Value = '18.95'
$ws.Range('E28').Value = '  +1.82%  '
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('D30').Value = '4.348'
$ws.Range('E30').Value = '  -2.02%  '
$ws.Range('D31').Value = '4.266'
$ws.Range('E31').Value = '  -2.29%  '
$ws.Range('D32').Value = '1.210'
$ws.Range('E32').Value = '  +0.53%  '
$ws.Range('D33').Value = '0.05324'
$ws.Range('E33').Value = '  -1.39%  '
$ws.Range('D34').Value = '1.904'
$ws.Range('E34').Value = '  -1.95%  '
$ws.Range('D35').Value = '1.177'
$ws.Range('E35').Value = '  -0.18%  '
$ws.Range('D36').Value = '0.7496'
$ws.Range('E36').Value = '  -3.51%  '
$ws.Range('E37').Value = '  +0.83%  '
$ws.Range('E38').Value = '  -0.12%  '
$ws.Range('D39').Value = '1.272.32'
$ws.Range('E39').Value = '  +0.20%  '
$ws.Range('D40').Value = '2.752'
$ws.Range('E40').Value = '  +0.27%  '
$ws.Range('D41').Value = '0.9052'
$ws.Range('E41').Value = '  -1.73%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').Value = '108.29'
$ws.Range('E42').Value = '  -4.86%  '
$ws.Range('D43').Value = '71.90'
$ws.Range('E43').Value = '  -3.69%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '5.983'
$ws.Range('E44').Value = '  -8.70%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').Value = '0.9997'
$ws.Range('E45').Value = '  -0.16%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.00000000130'
$ws.Range('E46').Value = '  +2.02%  '
$ws.Range('D47').Value = '2.024.85'
$ws.Range('E47').Value = '  -1.22%  '
$ws.Range('D48').Value = '1.797'
$ws.Range('E48').Value = '  -0.48%  '
$ws.Range('D49').Value = '0.5192'
$ws.Range('E49').Value = '  -0.53%  '
$ws.Range('D50').Value = '9.545'
$ws.Range('E50').Value = '  -0.10%  '
$ws.Range('D51').Value = '0.4337'
$ws.Range('E51').Value = '  -0.98%  '

# Restore default ("Normal") style on the cells we touched above so the
# saved workbook keeps the original (unstyled) look for these data cells.
foreach ($addr in $numericRiskCells) {
    $ws.Range($addr).Style = "Normal"
}
